$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
# Cells whose new value looks like a plain number (e.g. "93.70") must be
# forced to Text format first, otherwise Excel auto-converts them to a
# numeric value and silently drops significant trailing zeros
# (e.g. "93.70" -> 93.7), which would not match the source data.
$ws.Range("D2").Value = "27.699.99"
$ws.Range("D3").Value = "1.778.19"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.27"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4600"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07491"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.01"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.107"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9996"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.88"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.042"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.226"
$ws.Range("D16").Value = "1.773.96"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.70"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9991"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.09"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.793"
$ws.Range("D23").Value = "27.769.19"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.083"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.75"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.34"
$ws.Range("D28").Value = "1.976.31"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.171"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.59"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.100"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09228"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.677"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.554"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02295"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06108"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6323"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.972"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.184"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.396"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.801"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.29"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.730"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5892"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.51"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.952"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06950"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.41"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("E3").Value = "  +1.77%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").Value = "  +3.98%  "
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("E23").Value = "  +0.91%  "
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("E26").Value = "  +1.96%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("E29").Value = "  +5.05%  "
$ws.Range("E30").Value = "  +2.09%  "
$ws.Range("E31").Value = "  +3.16%  "
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("E34").Value = "  +1.69%  "
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("E40").Value = "  +1.63%  "
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("E43").Value = "  +0.82%  "
$ws.Range("E44").Value = "  +0.34%  "
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("E47").Value = "  +1.39%  "
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("E51").Value = "  +0.91%  "
